# Minor edits to proposal
$d = $word.ActiveDocument

# 1) Author line: fix the spelling "Hassaj" -> "Hasaj"
$d.Content.Find.Execute(
    "Hassaj",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Hasaj", 2) | Out-Null

# 2) Proposed Solution: expand the sentence with "(1) ... (2) ..." framing
$d.Content.Find.Execute(
    ": Provide Blue Bikes a list of bike stands that will benefit from interventions that will improve customer experience and potentially reduce costs and provide brief explanations of these solutions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Provide Blue Bikes with (1) a list of bike stands that will benefit from interventions aimed to both improve customer experience and potentially reduce costs and (2) provide brief explanations of these solutions.",
    2) | Out-Null

# 3) Make the "Methodology" heading word bold (keeping its existing italics)
$rngMethodology = $d.Paragraphs.Item(8).Range.Words.Item(1)
$rngMethodology.Font.Bold = 1
$rngMethodology.Font.BoldBi = 1

# 4) Predictive model bullet rewording
$d.Content.Find.Execute(
    "Build a predictive model: Given any time of day what are the amounts of bikes at a station ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Build a predictive model: Given any time of day how many bikes are at a station.",
    2) | Out-Null

# 5) "Given that a station..." bullet rewording
$d.Content.Find.Execute(
    "Given that a station X at time Y has no bikes, how long will it take for on bike to arrive. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Given that station X at time Y has no bikes, how long will it take for a bike to arrive. ",
    2) | Out-Null

# 6) "Provide intervention analysis on these stations" -> "...for these stations"
$d.Content.Find.Execute(
    "Provide intervention analysis on these stations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Provide intervention analysis for these stations",
    2) | Out-Null

# 7) "Refill the bike station more" -> "...more often"
$d.Content.Find.Execute(
    "Refill the bike station more",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Refill the bike station more often",
    2) | Out-Null

# 8) "Increase already in place infrastructure..." -> "Increase current infrastructure..."
$d.Content.Find.Execute(
    "Increase already in place infrastructure to restock bike stands more than currently done. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Increase current infrastructure to restock bike stands more than currently done. ",
    2) | Out-Null

# 9) "What are the top 15 single use stations." wording unchanged, but re-touch run
$d.Content.Find.Execute(
    "What are the top 15 single use stations. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "What are the top 15 single use stations. ",
    2) | Out-Null

# 10) "When are the most single uses occurring." wording unchanged, but re-touch run
$d.Content.Find.Execute(
    "When are the most single uses occurring. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "When are the most single uses occurring. ",
    2) | Out-Null
